$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the frequency response table (columns I and J)
$ws.Range("I1").Value = "f, kHz"
$ws.Range("I1").NumberFormat = "0.00"

$ws.Range("J1").Value = "V2, V"
$ws.Range("J1").NumberFormat = "0.0000"

# Frequency (kHz) values, column I, rows 2-32 (100..160 step 2)
$iVals = @(100,102,104,106,108,110,112,114,116,118,120,122,124,126,128,130,132,134,136,138,140,142,144,146,148,150,152,154,156,158,160)

# Corresponding V2 (V) readings, column J, rows 2-32
$jVals = @(0.16,0.16,0.16,0.16,0.16,0.16,0.17499999999999999,0.2,0.22500000000000001,0.25,0.3,0.4,0.5,0.75,1.21,1.5,0.8,0.48,0.33,0.22500000000000001,0.17499999999999999,0.15,0.11,0.1,0.074999999999999997,0.08,0.05,0.04,0.3,0.25,0.2)

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $iVals[$i]

    $jCell = $ws.Cells.Item($row, 10)
    $jCell.Value = $jVals[$i]
    $jCell.NumberFormat = "0.0000"
}

# Size column J to fit its (longer, 4-decimal) contents
$ws.Columns.Item(10).ColumnWidth = 9.7

# Match the recorded selection left by the author after the edit
$ws.Range("L24").Select() | Out-Null
